# ProjectManagement.xlsx edit script
# Commit message: "add baseline data and check timeline"
#
# Summary of changes:
#  - Timeline sheet: add two checkpoint rows (mail / all data done) with dates
#    in columns A/B (rows 15-16), formatted with the "d-mmm" number format.
#  - Cost sheet: update Treatment Units baseline (B1: 200 -> 100), which
#    ripples into the dependent cost formulas; add a new baseline row 12
#    (A4 / flyer, post card).
#  - Selections / active sheet updated to reflect where the author ended up
#    (Cost sheet active, Timeline selection at F13, Cost selection at B12).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Timeline
$ws2 = $wb.Worksheets.Item(2)   # Cost

# --- Timeline sheet: new checkpoint entries (rows 15 & 16, columns A/B) ---
$ws1.Select()

$ws1.Range("A15").Value = "mail"
$ws1.Range("B15").Value = 43182
$ws1.Range("B15").NumberFormat = "d-mmm"

$ws1.Range("A16").Value = "all data done"
$ws1.Range("B16").Value = 43202
$ws1.Range("B16").NumberFormat = "d-mmm"

# leave the selection on the Timeline sheet where the author left it
$ws1.Range("F13").Select()

# --- Cost sheet: baseline data updates ---
$ws2.Select()

# Treatment Units baseline changes from 200 to 100 (recalculates C3/C6/C8)
$ws2.Range("B1").Value = 100

# new row of baseline data (entered B12 before A12)
$ws2.Range("B12").Value = "flyer, post card"
$ws2.Range("A12").Value = "A4"

$ws2.Range("B12").Select()

# Cost sheet ends up being the active sheet
$ws2.Activate()
